# Update the "想去人数" (interest count) column F on the 展览 (sheet1)
# and 全部类型 (sheet4) sheets: each of these rows gains +1 attendee.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)   # 展览
$sheet4 = $wb.Worksheets.Item(4)   # 全部类型

# Row -> new value, for the "展览" sheet
$updates1 = @{
    5  = 1820
    9  = 2362
    16 = 313
    18 = 17
    20 = 196
    25 = 34
    26 = 1476
    29 = 229
}

foreach ($row in $updates1.Keys) {
    $sheet1.Cells.Item($row, 6).Value = $updates1[$row]
}

# Row -> new value, for the "全部类型" sheet (one extra row offsets indices by +1)
$updates4 = @{
    5  = 1820
    10 = 2362
    17 = 313
    19 = 17
    21 = 196
    26 = 34
    27 = 1476
    30 = 229
}

foreach ($row in $updates4.Keys) {
    $sheet4.Cells.Item($row, 6).Value = $updates4[$row]
}
